$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.453.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.34%  "

$ws.Range("D3").Value = "'3.369.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'590.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.13%  "

$ws.Range("D6").Value = "'187.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.20%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.37%  "

$ws.Range("D9").Value = "'0.183"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").Value = "'0.587"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("D11").Value = "'47.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.69%  "

$ws.Range("D12").Value = "'0.0000274"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.96%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "'3.915.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").Value = "'640.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.75%  "

$ws.Range("D15").Value = "'8.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.17%  "

$ws.Range("D16").Value = "'67.487.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.45%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.375.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.119"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("D19").Value = "'18.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "'11.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.71%  "

$ws.Range("D21").Value = "'0.910"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("D22").Value = "'18.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.53%  "

$ws.Range("D23").Value = "'5.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "'100.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").Value = "'4.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("E26").Value = "  +2.89%  "

$ws.Range("D27").Value = "'9.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("D28").Value = "'32.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.72%  "

$ws.Range("D29").Value = "'8.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").Value = "'6.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("D31").Value = "'612.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.80%  "

$ws.Range("D32").Value = "'3.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.58%  "

$ws.Range("D33").Value = "'11.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").Value = "'3.921.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.68%  "

$ws.Range("D35").Value = "'0.106"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.43%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").Value = "'55.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.15%  "

$ws.Range("D38").Value = "'2.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.05%  "

$ws.Range("E39").Value = "  +2.14%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'3.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'33.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("D42").Value = "'0.0₃0705"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.12%  "

$ws.Range("D43").Value = "'0.344"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("D44").Value = "'3.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.46%  "

$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'1.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.50%  "

$ws.Range("B50").Value = "CoreDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D50").Value = "'2.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -21.49%  "

$ws.Range("D51").Value = "'129.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.32%  "
